# Generate Report for Handoff
# Update the localization-status report so it reflects the newly generated
# handoff artifacts: new source GUID, new xliff content hashes, and the
# refreshed "Ready for handoff" / xliff-generation timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "50b5a56d-b964-4196-9f07-a2f86a41e16b"
$newGuid = "77ce5d41-ba5b-4fc6-97bf-6b987d6f2067"

$oldZhXlf = "$oldGuid.8feb56730504511871ca86bd31f05a0c150aea90.zh-cn.xlf"
$newZhXlf = "$newGuid.4ecf726fea8c58579321d739c736d4ed892de1b9.zh-cn.xlf"

$oldDeXlf = "$oldGuid.8feb56730504511871ca86bd31f05a0c150aea90.de-de.xlf"
$newDeXlf = "$newGuid.4ecf726fea8c58579321d739c736d4ed892de1b9.de-de.xlf"

$newHandoffDateTime = "2016-09-05 01:04:34"
$newZhGenerated = "2016-09-05 01:04:29"

# The hyperlinks in this workbook are EXTERNAL (they point out to GitHub);
# their target URL is unchanged by this edit, only the visible display text
# needs to move to the new filename. External hyperlinks can't be patched
# in place, so each is deleted and re-added against the same Address with
# the refreshed TextToDisplay.

function Set-HyperlinkDisplay {
    param($ws, $cellRef, $address, $display)
    $range = $ws.Range($cellRef)
    $range.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($range, $address, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $display) | Out-Null
}

$githubAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/42ef0280187a55dd1ee19a305e6dc919e94bb742/e2e/$oldGuid.md"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Sheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
Set-HyperlinkDisplay $wsOverview "B2" $githubAddress "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = $newHandoffDateTime

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Sheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
Set-HyperlinkDisplay $wsZhCn "A2" $githubAddress "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhXlf
$wsZhCn.Range("H2").Value = $newZhGenerated

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Sheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
Set-HyperlinkDisplay $wsDeDe "A2" $githubAddress "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeXlf
$wsDeDe.Range("H2").Value = $newHandoffDateTime
